# Regenerate the "K" (strikeouts) column (G) values for the 2022 save_data
# sheet. The previous data used a different strikeout counting method
# ("Strike#"); this updates column G to the recomputed "K" values for each
# game row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 2
    6  = 1
    7  = 1
    8  = 1
    9  = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 0
    14 = 2
    16 = 0
    17 = 2
    18 = 0
    19 = 0
    20 = 0
    21 = 2
    22 = 2
    23 = 2
    24 = 0
    25 = 3
    26 = 2
    27 = 0
    28 = 0
    29 = 0
    30 = 1
    31 = 3
    32 = 3
    33 = 1
    36 = 1
    37 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
